# "Ajout Buzzer + main testé"
#
# The "Detecteur de porte" objective (row 3 of the Objectifs table) moves
# from "à tester" to "validé" in the Progres column (G).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objectifs")

$ws.Range("G3").Value = "validé"

# Leave the selection where the author left it when they saved.
[void]$ws.Range("G11").Select()
